$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply corrected values and matching fill styles cell by cell.
# Style legend (fill color only; border/alignment unchanged):
#   style 2 -> orange (FFC966), style 3 -> green (6DC066),
#   style 4 -> red (FF6666),   style 5 -> purple (8067A2)
# Template cells below keep a stable, never-edited copy of each
# fill style so Copy+PasteSpecial reuses the workbook's existing
# style (xf) index instead of synthesizing a new one per cell.
$styleTemplate2 = "B4"
$styleTemplate3 = "D4"
$styleTemplate4 = "G4"
$styleTemplate5 = "K33"

$ws.Range("E4").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("E4").PasteSpecial(-4122)

$ws.Range("H4").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("H4").PasteSpecial(-4122)

$ws.Range("P4").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("P4").PasteSpecial(-4122)

$ws.Range("C5").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("C5").PasteSpecial(-4122)

$ws.Range("D5").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("D5").PasteSpecial(-4122)

$ws.Range("E5").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("E5").PasteSpecial(-4122)

$ws.Range("F5").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("F5").PasteSpecial(-4122)

$ws.Range("H5").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("H5").PasteSpecial(-4122)

$ws.Range("I5").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("I5").PasteSpecial(-4122)

$ws.Range("M5").Value = 2
$ws.Range($styleTemplate3).Copy()
$ws.Range("M5").PasteSpecial(-4122)

$ws.Range("P5").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("P5").PasteSpecial(-4122)

$ws.Range("C6").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("C6").PasteSpecial(-4122)

$ws.Range("D6").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("D6").PasteSpecial(-4122)

$ws.Range("F6").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("F6").PasteSpecial(-4122)

$ws.Range("G6").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("G6").PasteSpecial(-4122)

$ws.Range("H6").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("H6").PasteSpecial(-4122)

$ws.Range("J6").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("J6").PasteSpecial(-4122)

$ws.Range("B7").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("B7").PasteSpecial(-4122)

$ws.Range("C7").Value = 2
$ws.Range($styleTemplate3).Copy()
$ws.Range("C7").PasteSpecial(-4122)

$ws.Range("E7").Value = 2
$ws.Range($styleTemplate3).Copy()
$ws.Range("E7").PasteSpecial(-4122)

$ws.Range("F7").Value = 2
$ws.Range($styleTemplate3).Copy()
$ws.Range("F7").PasteSpecial(-4122)

$ws.Range("H7").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("H7").PasteSpecial(-4122)

$ws.Range("I7").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("I7").PasteSpecial(-4122)

$ws.Range("J7").Value = 2
$ws.Range($styleTemplate3).Copy()
$ws.Range("J7").PasteSpecial(-4122)

$ws.Range("G8").Value = 2
$ws.Range($styleTemplate3).Copy()
$ws.Range("G8").PasteSpecial(-4122)

$ws.Range("H8").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("H8").PasteSpecial(-4122)

$ws.Range("I8").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("I8").PasteSpecial(-4122)

$ws.Range("J8").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("J8").PasteSpecial(-4122)

$ws.Range("L8").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("L8").PasteSpecial(-4122)

$ws.Range("P8").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("P8").PasteSpecial(-4122)

$ws.Range("B9").Value = 2
$ws.Range($styleTemplate3).Copy()
$ws.Range("B9").PasteSpecial(-4122)

$ws.Range("C9").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("C9").PasteSpecial(-4122)

$ws.Range("E9").Value = 2
$ws.Range($styleTemplate3).Copy()
$ws.Range("E9").PasteSpecial(-4122)

$ws.Range("F9").Value = 2
$ws.Range($styleTemplate3).Copy()
$ws.Range("F9").PasteSpecial(-4122)

$ws.Range("H9").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("H9").PasteSpecial(-4122)

$ws.Range("I9").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("I9").PasteSpecial(-4122)

$ws.Range("J9").Value = 2
$ws.Range($styleTemplate3).Copy()
$ws.Range("J9").PasteSpecial(-4122)

$ws.Range("P9").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("P9").PasteSpecial(-4122)

$ws.Range("G10").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("G10").PasteSpecial(-4122)

$ws.Range("I10").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("I10").PasteSpecial(-4122)

$ws.Range("J10").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("J10").PasteSpecial(-4122)

$ws.Range("K10").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("K10").PasteSpecial(-4122)

$ws.Range("B11").Value = 2
$ws.Range($styleTemplate3).Copy()
$ws.Range("B11").PasteSpecial(-4122)

$ws.Range("C11").Value = 2
$ws.Range($styleTemplate3).Copy()
$ws.Range("C11").PasteSpecial(-4122)

$ws.Range("H11").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("H11").PasteSpecial(-4122)

$ws.Range("J11").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("J11").PasteSpecial(-4122)

$ws.Range("K11").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("K11").PasteSpecial(-4122)

$ws.Range("P11").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("P11").PasteSpecial(-4122)

$ws.Range("F12").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("F12").PasteSpecial(-4122)

$ws.Range("I12").Value = 2
$ws.Range($styleTemplate3).Copy()
$ws.Range("I12").PasteSpecial(-4122)

$ws.Range("J12").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("J12").PasteSpecial(-4122)

$ws.Range("P12").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("P12").PasteSpecial(-4122)

$ws.Range("G13").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("G13").PasteSpecial(-4122)

$ws.Range("I13").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("I13").PasteSpecial(-4122)

$ws.Range("P13").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("P13").PasteSpecial(-4122)

$ws.Range("C18").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("C18").PasteSpecial(-4122)

$ws.Range("D18").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("D18").PasteSpecial(-4122)

$ws.Range("F18").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("F18").PasteSpecial(-4122)

$ws.Range("I18").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("I18").PasteSpecial(-4122)

$ws.Range("B19").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("B19").PasteSpecial(-4122)

$ws.Range("D19").Value = 2
$ws.Range($styleTemplate3).Copy()
$ws.Range("D19").PasteSpecial(-4122)

$ws.Range("G19").Value = 2
$ws.Range($styleTemplate3).Copy()
$ws.Range("G19").PasteSpecial(-4122)

$ws.Range("B20").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("B20").PasteSpecial(-4122)

$ws.Range("D20").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("D20").PasteSpecial(-4122)

$ws.Range("F20").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("F20").PasteSpecial(-4122)

$ws.Range("G20").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("G20").PasteSpecial(-4122)

$ws.Range("I20").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("I20").PasteSpecial(-4122)

$ws.Range("B21").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("B21").PasteSpecial(-4122)

$ws.Range("C21").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("C21").PasteSpecial(-4122)

$ws.Range("D21").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("D21").PasteSpecial(-4122)

$ws.Range("E21").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("E21").PasteSpecial(-4122)

$ws.Range("G21").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("G21").PasteSpecial(-4122)

$ws.Range("B22").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("B22").PasteSpecial(-4122)

$ws.Range("C22").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("D22").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("D22").PasteSpecial(-4122)

$ws.Range("F22").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("F22").PasteSpecial(-4122)

$ws.Range("H22").Value = 2
$ws.Range($styleTemplate3).Copy()
$ws.Range("H22").PasteSpecial(-4122)

$ws.Range("I22").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("I22").PasteSpecial(-4122)

$ws.Range("B23").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("B23").PasteSpecial(-4122)

$ws.Range("C23").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("C23").PasteSpecial(-4122)

$ws.Range("F23").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("F23").PasteSpecial(-4122)

$ws.Range("I23").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("I23").PasteSpecial(-4122)

$ws.Range("C24").Value = 2
$ws.Range($styleTemplate3).Copy()
$ws.Range("C24").PasteSpecial(-4122)

$ws.Range("D24").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("D24").PasteSpecial(-4122)

$ws.Range("E24").Value = 2
$ws.Range($styleTemplate3).Copy()
$ws.Range("E24").PasteSpecial(-4122)

$ws.Range("F24").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("F24").PasteSpecial(-4122)

$ws.Range("H24").Value = 2
$ws.Range($styleTemplate3).Copy()
$ws.Range("H24").PasteSpecial(-4122)

$ws.Range("C25").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("C25").PasteSpecial(-4122)

$ws.Range("B26").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("B26").PasteSpecial(-4122)

$ws.Range("G26").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("G26").PasteSpecial(-4122)

$ws.Range("G27").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("G27").PasteSpecial(-4122)

$ws.Range("D32").Value = 2
$ws.Range($styleTemplate3).Copy()
$ws.Range("D32").PasteSpecial(-4122)

$ws.Range("F32").Value = 2
$ws.Range($styleTemplate3).Copy()
$ws.Range("F32").PasteSpecial(-4122)

$ws.Range("G32").Value = 2
$ws.Range($styleTemplate3).Copy()
$ws.Range("G32").PasteSpecial(-4122)

$ws.Range("B33").Value = 3
$ws.Range($styleTemplate5).Copy()
$ws.Range("B33").PasteSpecial(-4122)

$ws.Range("D33").Value = 2
$ws.Range($styleTemplate3).Copy()
$ws.Range("D33").PasteSpecial(-4122)

$ws.Range("J33").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("J33").PasteSpecial(-4122)

$ws.Range("D34").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("D34").PasteSpecial(-4122)

$ws.Range("C35").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("C35").PasteSpecial(-4122)

$ws.Range("D35").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("D35").PasteSpecial(-4122)

$ws.Range("E35").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("E35").PasteSpecial(-4122)

$ws.Range("G35").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("G35").PasteSpecial(-4122)

$ws.Range("H35").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("H35").PasteSpecial(-4122)

$ws.Range("B36").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("B36").PasteSpecial(-4122)

$ws.Range("C36").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("C36").PasteSpecial(-4122)

$ws.Range("D36").Value = 2
$ws.Range($styleTemplate3).Copy()
$ws.Range("D36").PasteSpecial(-4122)

$ws.Range("E36").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("E36").PasteSpecial(-4122)

$ws.Range("G36").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("G36").PasteSpecial(-4122)

$ws.Range("H36").Value = 2
$ws.Range($styleTemplate3).Copy()
$ws.Range("H36").PasteSpecial(-4122)

$ws.Range("C37").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("C37").PasteSpecial(-4122)

$ws.Range("D37").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("D37").PasteSpecial(-4122)

$ws.Range("G37").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("G37").PasteSpecial(-4122)

$ws.Range("C38").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("C38").PasteSpecial(-4122)

$ws.Range("F38").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("F38").PasteSpecial(-4122)

$ws.Range("H38").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("H38").PasteSpecial(-4122)

$ws.Range("J38").Value = 3
$ws.Range($styleTemplate5).Copy()
$ws.Range("J38").PasteSpecial(-4122)

$ws.Range("B39").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("B39").PasteSpecial(-4122)

$ws.Range("C39").Value = 0
$ws.Range($styleTemplate4).Copy()
$ws.Range("C39").PasteSpecial(-4122)

$ws.Range("F39").Value = 2
$ws.Range($styleTemplate3).Copy()
$ws.Range("F39").PasteSpecial(-4122)

$ws.Range("H39").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("H39").PasteSpecial(-4122)

$ws.Range("B40").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("B40").PasteSpecial(-4122)

$ws.Range("H40").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("H40").PasteSpecial(-4122)

$ws.Range("I40").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("I40").PasteSpecial(-4122)

$ws.Range("B41").Value = 2
$ws.Range($styleTemplate3).Copy()
$ws.Range("B41").PasteSpecial(-4122)

$ws.Range("H41").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("H41").PasteSpecial(-4122)

$ws.Range("I41").Value = 1
$ws.Range($styleTemplate2).Copy()
$ws.Range("I41").PasteSpecial(-4122)

$excel.CutCopyMode = 0